$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$errorText = @'
no such element: Unable to locate element: {"method":"xpath","selector":".//*[@id='wsGrid3']/div[1]/div[1]/div[2]/div/div[1]/div/div[8]/div/a/span"}
  (Session info: chrome=56.0.2924.87)
  (Driver info: chromedriver=2.25.426923 (0390b88869384d6eb0d5d09729679f934aab9eed),platform=Windows NT 10.0.14393 x86_64) (WARNING: The server did not provide any stacktrace information)
Command duration or timeout: 10.35 seconds
For documentation on this error, please visit: http://seleniumhq.org/exceptions/no_such_element.html
Build info: version: '3.0.1', revision: '1969d75', time: '2016-10-18 09:49:13 -0700'
System info: host: 'MQCSERVER', ip: '172.16.0.6', os.name: 'Windows 10', os.arch: 'amd64', os.version: '10.0', java.version: '1.8.0_121'
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities [{applicationCacheEnabled=false, rotatable=false, mobileEmulationEnabled=false, networkConnectionEnabled=false, chrome={chromedriverVersion=2.25.426923 (0390b88869384d6eb0d5d09729679f934aab9eed), userDataDir=C:\Users\admin\AppData\Local\Temp\scoped_dir4028_26424}, takesHeapSnapshot=true, pageLoadStrategy=normal, databaseEnabled=false, handlesAlerts=true, hasTouchScreen=false, version=56.0.2924.87, platform=XP, browserConnectionEnabled=false, nativeEvents=true, acceptSslCerts=true, locationContextEnabled=true, webStorageEnabled=true, browserName=chrome, takesScreenshot=true, javascriptEnabled=true, cssSelectorsEnabled=true}]
Session ID: 95b176ce5c6d8567fb219bf37ddbfd7f
*** Element info: {Using=xpath, value=.//*[@id='wsGrid3']/div[1]/div[1]/div[2]/div/div[1]/div/div[8]/div/a/span}
'@

# Row 9: duplicate of the TC001 / created-and-approved Pass row
$ws.Range("A9").Value = "TC001"
$ws.Range("B9").Value = "Creating the TEAM Workspace"
$ws.Range("C9").Value = "TEAM Workspace Should be created successfully and approved"
$ws.Range("D9").Value = "TEAM workspace is created successfully and approved"
$ws.Range("E9").Value = "Pass"

# Row 10: duplicate of the TC001 / created-and-approved Pass row
$ws.Range("A10").Value = "TC001"
$ws.Range("B10").Value = "Creating the TEAM Workspace"
$ws.Range("C10").Value = "TEAM Workspace Should be created successfully and approved"
$ws.Range("D10").Value = "TEAM workspace is created successfully and approved"
$ws.Range("E10").Value = "Pass"

# Row 11: new Fail row with the selenium no-such-element stack trace
$ws.Range("A11").Value = "TC001"
$ws.Range("B11").Value = "Creating the TEAM Workspace"
$ws.Range("C11").Value = "TEAM Workspace Should be created successfully"
$ws.Range("D11").Value = $errorText
$ws.Range("E11").Value = "Fail"
